# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on a handful of Leve rows
# across the per-class Sheets, mirroring the latest Universalis snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 205.4
$ws.Range("I12").Value = 98.3
$ws.Range("K12").Value = 98.3
$ws.Range("M12").Value = 71.7

$ws.Range("H32").Value = 882.2222
$ws.Range("I32").Value = 650
$ws.Range("J32").Value = 928.6667
$ws.Range("K32").Value = 650
$ws.Range("L32").Value = 928.6667
$ws.Range("M32").Value = -324
$ws.Range("N32").Value = -1580.6667

$ws.Range("H86").Value = 58826492
$ws.Range("I86").Value = 2215.875
$ws.Range("J86").Value = 111114740
$ws.Range("K86").Value = 2215.875
$ws.Range("L86").Value = 111114740
$ws.Range("M86").Value = -1092.875
$ws.Range("N86").Value = -111116986

$ws.Range("H89").Value = 58826492
$ws.Range("I89").Value = 2215.875
$ws.Range("J89").Value = 111114740
$ws.Range("K89").Value = 11079.375
$ws.Range("L89").Value = 555573700
$ws.Range("M89").Value = -5463.375
$ws.Range("N89").Value = -555584932

$ws.Range("H116").Value = 5687.375
$ws.Range("I116").Value = 6874.75
$ws.Range("J116").Value = 4500
$ws.Range("K116").Value = 6874.75
$ws.Range("L116").Value = 4500
$ws.Range("M116").Value = -3432.75
$ws.Range("N116").Value = -11384

$ws.Range("H128").Value = 35708.57
$ws.Range("J128").Value = 35708.57
$ws.Range("L128").Value = 35708.57
$ws.Range("N128").Value = -45668.57

$ws.Range("H132").Value = 1669.4318
$ws.Range("I132").Value = 1669.4318
$ws.Range("K132").Value = 5008.2954
$ws.Range("M132").Value = -2478.2954

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10725.549
$ws.Range("I32").Value = 8646.413
$ws.Range("K32").Value = 8646.413
$ws.Range("M32").Value = -8359.413

$ws.Range("H45").Value = 961.8823
$ws.Range("I45").Value = 809
$ws.Range("J45").Value = 1097.7778
$ws.Range("K45").Value = 809
$ws.Range("L45").Value = 1097.7778
$ws.Range("M45").Value = -432
$ws.Range("N45").Value = -1851.7778

$ws.Range("H61").Value = 2203.8667
$ws.Range("I61").Value = 2074.5715
$ws.Range("J61").Value = 4014
$ws.Range("K61").Value = 2074.5715
$ws.Range("L61").Value = 4014
$ws.Range("M61").Value = -1862.5715
$ws.Range("N61").Value = -4438

$ws.Range("H74").Value = 1793.1428
$ws.Range("I74").Value = 1777.2307
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1777.2307
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -903.2307000000001
$ws.Range("N74").Value = -3748

$ws.Range("H77").Value = 1793.1428
$ws.Range("I77").Value = 1777.2307
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 8886.1535
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -4518.1535
$ws.Range("N77").Value = -18736

$ws.Range("H132").Value = 2423.9143
$ws.Range("I132").Value = 1880.48
$ws.Range("K132").Value = 5641.440000000001
$ws.Range("M132").Value = -3111.440000000001

$ws.Range("H136").Value = 2203.8667
$ws.Range("I136").Value = 2074.5715
$ws.Range("J136").Value = 4014
$ws.Range("K136").Value = 6223.7145
$ws.Range("L136").Value = 12042
$ws.Range("M136").Value = -3673.7145
$ws.Range("N136").Value = -17142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 18888.95
$ws.Range("I134").Value = 1597.8605
$ws.Range("K134").Value = 4793.5815
$ws.Range("M134").Value = -2258.5815

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1196.75
$ws.Range("I22").Value = 1637.2858
$ws.Range("J22").Value = 580
$ws.Range("K22").Value = 1637.2858
$ws.Range("L22").Value = 580
$ws.Range("M22").Value = -1287.2858
$ws.Range("N22").Value = -1280

$ws.Range("H86").Value = 12100
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377

$ws.Range("H89").Value = 12100
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884

$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = 3600
$ws.Range("N122").Value = -8500

$ws.Range("H132").Value = 2056.5264
$ws.Range("I132").Value = 1289.3
$ws.Range("J132").Value = 2909
$ws.Range("K132").Value = 3867.9
$ws.Range("L132").Value = 8727
$ws.Range("M132").Value = -1337.9
$ws.Range("N132").Value = -13787

$ws.Range("H134").Value = 12821731
$ws.Range("I134").Value = 1126.1471
$ws.Range("K134").Value = 3378.4413
$ws.Range("M134").Value = -843.4412999999995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 869.95557
$ws.Range("I122").Value = 346.66666
$ws.Range("J122").Value = 907.3333
$ws.Range("K122").Value = 3119.99994
$ws.Range("L122").Value = 8165.9997
$ws.Range("M122").Value = -669.9999399999997
$ws.Range("N122").Value = -13065.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4023.0256
$ws.Range("I70").Value = 4005.9033
$ws.Range("K70").Value = 4005.9033
$ws.Range("M70").Value = -3735.9033

$ws.Range("H73").Value = 4023.0256
$ws.Range("I73").Value = 4005.9033
$ws.Range("K73").Value = 4005.9033
$ws.Range("M73").Value = -3069.9033

$ws.Range("H123").Value = 41941.066
$ws.Range("J123").Value = 41941.066
$ws.Range("L123").Value = 41941.066
$ws.Range("N123").Value = -46841.066

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1187.125
$ws.Range("I22").Value = 600.5
$ws.Range("J22").Value = 1270.9286
$ws.Range("K22").Value = 600.5
$ws.Range("L22").Value = 1270.9286
$ws.Range("M22").Value = -305.5
$ws.Range("N22").Value = -1860.9286

$ws.Range("H27").Value = 1187.125
$ws.Range("I27").Value = 600.5
$ws.Range("J27").Value = 1270.9286
$ws.Range("K27").Value = 600.5
$ws.Range("L27").Value = 1270.9286
$ws.Range("M27").Value = -493.5
$ws.Range("N27").Value = -1484.9286

$ws.Range("H40").Value = 2526566.2
$ws.Range("I40").Value = 3368323.2
$ws.Range("K40").Value = 3368323.2
$ws.Range("M40").Value = -3368187.2

$ws.Range("H61").Value = 2313.3333
$ws.Range("I61").Value = 3296.6667
$ws.Range("K61").Value = 3296.6667
$ws.Range("M61").Value = -3094.6667

$ws.Range("H113").Value = 2313.3333
$ws.Range("I113").Value = 3296.6667
$ws.Range("K113").Value = 3296.6667
$ws.Range("M113").Value = -1126.6667

$ws.Range("H132").Value = 3244.389
$ws.Range("I132").Value = 2760
$ws.Range("K132").Value = 8280
$ws.Range("M132").Value = -5750

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1223.2812
$ws.Range("I132").Value = 1256.2778
$ws.Range("J132").Value = 1180.8572
$ws.Range("K132").Value = 3768.8334
$ws.Range("L132").Value = 3542.5716
$ws.Range("M132").Value = -1238.8334
$ws.Range("N132").Value = -8602.571599999999

$ws.Range("H136").Value = 1188.0555
$ws.Range("I136").Value = 1081.1538
$ws.Range("K136").Value = 3243.4614
$ws.Range("M136").Value = -693.4614000000001
